# Update the attendee/view counts in column F on the "展览" (sheet1) and
# "全部类型" (sheet4) worksheets to match the refreshed data pull.
# Row numbering differs by one between the two sheets from row 7 onward
# because "全部类型" has one extra data row ("展览" row 7 == "全部类型" row 8),
# so each sheet gets its own explicit row->value map.

$wb = $excel.ActiveWorkbook

$updates1 = @{
    3  = 2220
    4  = 93
    5  = 13306
    6  = 76
    7  = 119
    9  = 485
    10 = 1193
    11 = 996
    13 = 14456
    15 = 173
    17 = 43
    20 = 6
    21 = 41
    22 = 1105
    25 = 5503
    27 = 709
    28 = 348
    29 = 27
    30 = 104
}

$updates4 = @{
    3  = 2220
    4  = 93
    5  = 13306
    6  = 76
    8  = 119
    10 = 485
    11 = 1193
    12 = 996
    14 = 14456
    16 = 173
    18 = 43
    21 = 6
    22 = 41
    23 = 1105
    26 = 5503
    28 = 709
    29 = 348
    30 = 27
    31 = 104
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
